# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
# Both sheets hold identical event data, so the same row/value updates apply to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    6  = 623
    10 = 400
    17 = 1069
    18 = 1433
    22 = 89
    28 = 294
    29 = 1659
    33 = 607
    35 = 3869
    37 = 450
    39 = 978
    40 = 84
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
